# Apply the refreshed cryptocurrency price / volume(1h) data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values are plain decimals (e.g. 210.62) which Excel would
# otherwise auto-convert to a number. A leading apostrophe (quote-prefix)
# forces the cell to stay text, matching the original inline-string cells.

$ws.Range('D2').Value = "26.518.48"
$ws.Range('E2').Value = "  -2.41%  "
$ws.Range('D3').Value = "1.581.61"
$ws.Range('E3').Value = "  -3.11%  "
$ws.Range('E4').Value = "  +0.28%  "
$ws.Range('D5').Value = "`'210.62"
$ws.Range('E5').Value = "  -2.74%  "
$ws.Range('D6').Value = "`'0.506"
$ws.Range('E6').Value = "  -2.08%  "
$ws.Range('E7').Value = "  +0.29%  "
$ws.Range('E8').Value = "  -2.22%  "
$ws.Range('D9').Value = "`'0.0617"
$ws.Range('E10').Value = "  -3.94%  "
$ws.Range('E11').Value = "  -2.11%  "
$ws.Range('D12').Value = "1.801.99"
$ws.Range('E12').Value = "  -3.03%  "
$ws.Range('D13').Value = "1.591.36"
$ws.Range('E13').Value = "  -1.99%  "
$ws.Range('E14').Value = "  -1.70%  "
$ws.Range('D15').Value = "`'0.527"
$ws.Range('E15').Value = "  -2.81%  "
$ws.Range('D16').Value = "`'63.74"
$ws.Range('D17').Value = "26.545.66"
$ws.Range('E17').Value = "  -2.06%  "
$ws.Range('E18').Value = "  -0.80%  "
$ws.Range('E19').Value = "  +0.15%  "
$ws.Range('D20').Value = "`'207.94"
$ws.Range('E20').Value = "  -3.11%  "
$ws.Range('D21').Value = "`'6.67"
$ws.Range('E21').Value = "  -3.37%  "
$ws.Range('E22').Value = "  -3.38%  "
$ws.Range('D23').Value = "`'2.37"
$ws.Range('E23').Value = "  -5.18%  "
$ws.Range('D24').Value = "`'8.88"
$ws.Range('E24').Value = "  -2.02%  "
$ws.Range('D25').Value = "`'146.18"
$ws.Range('E25').Value = "  -1.53%  "
$ws.Range('E26').Value = "  +0.27%  "
$ws.Range('D27').Value = "`'7.42"
$ws.Range('E27').Value = "  +1.80%  "
$ws.Range('E28').Value = "  -4.46%  "
$ws.Range('D29').Value = "`'15.23"
$ws.Range('E29').Value = "  -2.21%  "
$ws.Range('E30').Value = "  -1.03%  "
$ws.Range('E31').Value = "  -2.52%  "
$ws.Range('D32').Value = "`'3.25"
$ws.Range('E32').Value = "  -3.77%  "
$ws.Range('D33').Value = "`'0.656"
$ws.Range('E33').Value = "  +22.19%  "
$ws.Range('D34').Value = "`'2.94"
$ws.Range('D35').Value = "1.305.53"
$ws.Range('E35').Value = "  -0.94%  "
$ws.Range('D36').Value = "`'1.50"
$ws.Range('E36').Value = "  -3.85%  "
$ws.Range('D37').Value = "`'2.43"
$ws.Range('E37').Value = "  -0.81%  "
$ws.Range('D38').Value = "`'0.0173"
$ws.Range('E38').Value = "  -0.91%  "
$ws.Range('D39').Value = "`'0.819"
$ws.Range('E39').Value = "  -2.98%  "
$ws.Range('E40').Value = "  +0.23%  "
$ws.Range('D41').Value = "`'0.783"
$ws.Range('E41').Value = "  -2.80%  "
$ws.Range('E42').Value = "  +0.82%  "
$ws.Range('E43').Value = "  -4.63%  "
$ws.Range('D44').Value = "`'62.68"
$ws.Range('D45').Value = "1.715.41"
$ws.Range('E45').Value = "  -2.80%  "
$ws.Range('E46').Value = "  -2.09%  "
$ws.Range('E47').Value = "  +0.45%  "
$ws.Range('D48').Value = "`'0.832"
$ws.Range('E48').Value = "  +5.76%  "
$ws.Range('E49').Value = "  -1.82%  "
$ws.Range('D50').Value = "`'0.0978"
$ws.Range('E50').Value = "  +3.30%  "
$ws.Range('B51').Value = "EnergySwap"
$ws.Range('C51').Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('D51').Value = "`'7.49"
$ws.Range('E51').Value = "  -0.78%  "
